$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename departments from underscore_separated to space separated labels.
# Update in an order that makes the new shared-string entries land at the
# end of the table in the same sequence as the target workbook:
#   madre de dios, san martin, la libertad
$ws.Range("B19").Value = "madre de dios"
$ws.Range("B24").Value = "san martin"
$ws.Range("B15").Value = "la libertad"

# Leave the final selection on B15, matching the saved cursor position.
$ws.Range("B15").Select()
